# Applies the cryptos list price/volume updates described in the commit message
# (prices in column D, 1h volume deltas in column E).
#
# Values are assigned with a leading apostrophe so Excel always stores them
# as literal text (matching the workbook's existing inlineStr cells) even
# when the text happens to look like a number (e.g. "571.18"); the style is
# then reset to Normal so the quote-prefix flag does not linger on the cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'60.980.43"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  +0.39%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'3.384.78"
$ws.Range("D3").Style = "Normal"
$ws.Range("D4").Value = "'0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'  +0.00%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'571.18"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  -0.05%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'141.87"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  +0.12%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("E7").Value = "'  +0.00%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("E8").Value = "'  +0.17%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("E9").Value = "'  +1.84%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("E10").Value = "'  -0.55%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("E11").Value = "'  -1.55%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'3.964.27"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  +0.07%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("E13").Value = "'  +1.74%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("E14").Value = "'  -0.43%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("E15").Value = "'  +0.68%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'3.386.98"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  +0.12%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'61.079.10"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  +0.36%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'6.09"
$ws.Range("D18").Style = "Normal"
$ws.Range("D19").Value = "'13.62"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  -3.46%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'8.90"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  -0.57%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'383.65"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  -1.17%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("E22").Value = "'  +2.66%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'0.553"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  -1.42%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'0.999"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  +0.26%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D26").Value = "'3.525.32"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  +0.06%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("E27").Value = "'  +2.45%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("E28").Value = "'  +0.07%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("E29").Value = "'  -1.61%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("E30").Value = "'  -1.56%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("E31").Value = "'  -1.09%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("E33").Value = "'  -5.14%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = "'23.19"
$ws.Range("D34").Style = "Normal"
$ws.Range("D35").Value = "'6.94"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  +0.23%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = "'166.45"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  -0.24%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = "'3.416.70"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  +0.20%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = "'4.98"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  -1.43%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("E39").Value = "'  -3.06%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("E40").Value = "'  -1.61%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'26.85"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  -0.19%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("E42").Value = "'  +0.05%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("E43").Value = "'  -0.55%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("E44").Value = "'  -1.97%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("E45").Value = "'  -1.72%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("E46").Value = "'  -0.17%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'2.445.75"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  -4.35%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'22.88"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  -0.01%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("E49").Value = "'  -2.11%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'2.16"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  +10.78%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("E51").Value = "'  +0.59%  "
$ws.Range("E51").Style = "Normal"

Write-Host "Applied cryptos list updates to 69 cells"
